$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values look numeric)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.561.52'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '1.869.53'
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("D4").Value = '1.000'

$ws.Range("D5").Value = '247.54'
$ws.Range("E5").Value = '  +1.23%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '0.4730'
$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("E8").Value = '  +1.45%  '

$ws.Range("D9").Value = '0.06469'
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").Value = '21.99'
$ws.Range("E10").Value = '  +4.65%  '

$ws.Range("D11").Value = '0.07703'
$ws.Range("E11").Value = '  -0.85%  '

$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '96.41'
$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7362'
$ws.Range("E13").Value = '  +3.55%  '

$ws.Range("D14").Value = '1.867.58'
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").Value = '5.137'
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '271.92'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").Value = '30.533.80'
$ws.Range("E17").Value = '  +0.80%  '

$ws.Range("D18").Value = '13.28'
$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").Value = '0.000007493'
$ws.Range("E20").Value = '  -0.39%  '

$ws.Range("D21").Value = '2.115.37'
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").Value = '5.235'
$ws.Range("E23").Value = '  +0.43%  '

$ws.Range("D24").Value = '6.162'
$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("D25").Value = '9.192'
$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("D26").Value = '163.74'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("D27").Value = '18.72'
$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("D28").Value = '1.904'
$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("D29").Value = '0.09958'
$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("D30").Value = '1.344'
$ws.Range("E30").Value = '  -2.67%  '

$ws.Range("D31").Value = '1.508'
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").Value = '4.271'
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("D33").Value = '4.098'
$ws.Range("E33").Value = '  +1.98%  '

$ws.Range("D34").Value = '0.04779'
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").Value = '0.6946'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.718'

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01846'
$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.748'
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.175'
$ws.Range("E40").Value = '  -2.21%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '72.89'
$ws.Range("E41").Value = '  +4.03%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '1.962'
$ws.Range("E42").Value = '  +2.72%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4169'
$ws.Range("E43").Value = '  +1.74%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8335'
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '101.17'
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.349'
$ws.Range("E47").Value = '  +1.25%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '35.43'
$ws.Range("E48").Value = '  +0.56%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '6.944'
$ws.Range("E49").Value = '  -2.00%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '915.45'
$ws.Range("E50").Value = '  -0.30%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05646'
$ws.Range("E51").Value = '  +1.34%  '

